$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update "Förändrad" (changed) date in column C for all existing data rows (2-475)
#    from 2023-09-21 (45190) to 2023-09-23 (45192)
$ws.Range("C2:C475").Value = 45192

# 2. Give row 475 the same explicit row height as all other data rows (15pt, custom)
$ws.Rows.Item(475).RowHeight = 15

# 3. Append a new record as row 476
$ws.Cells.Item(476, 1).Value = "A 44897-2023"
$ws.Cells.Item(476, 2).Value = 45190
$ws.Cells.Item(476, 3).Value = 45192
$ws.Cells.Item(476, 4).Value = "VÄSTRA GÖTALANDS LÄN"
$ws.Cells.Item(476, 5).Value = "MARK"
$ws.Cells.Item(476, 7).Value = 1
$ws.Cells.Item(476, 8).Value = 0
$ws.Cells.Item(476, 9).Value = 0
$ws.Cells.Item(476, 10).Value = 0
$ws.Cells.Item(476, 11).Value = 0
$ws.Cells.Item(476, 12).Value = 0
$ws.Cells.Item(476, 13).Value = 0
$ws.Cells.Item(476, 14).Value = 0
$ws.Cells.Item(476, 15).Value = 0
$ws.Cells.Item(476, 16).Value = 0
$ws.Cells.Item(476, 17).Value = 0

# Match number formatting used by the other rows (date columns, wrap-text in R)
$ws.Range("B476:C476").NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(476, 18).WrapText = $true
